$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-08 Sunday" "2023-10-09 Monday"

Replace-Text "67×84=" "19×17="
Replace-Text "43×73=" "60×99="
Replace-Text "71×21=" "20×65="
Replace-Text "55×75=" "60×91="
Replace-Text "64×74=" "60×61="

Replace-Text "81×76=" "53×35="
Replace-Text "53×97=" "54×84="
Replace-Text "21×13=" "75×86="
Replace-Text "54×56=" "23×38="
Replace-Text "35×36=" "83×37="

Replace-Text "23×15=" "83×21="
Replace-Text "61×22=" "11×15="
Replace-Text "74×11=" "49×90="
Replace-Text "49×28=" "57×63="
Replace-Text "66×55=" "58×86="

Replace-Text "56×14=" "21×21="
Replace-Text "42×95=" "90×32="
Replace-Text "47×82=" "75×88="
Replace-Text "90×23=" "32×28="
Replace-Text "37×29=" "57×58="

Replace-Text "36×13=" "39×69="
Replace-Text "43×71=" "32×39="
Replace-Text "77×48=" "39×38="
Replace-Text "89×99=" "20×65="
Replace-Text "97×17=" "88×67="
